# "Colocando header nos gráficos" (Adding headers for the charts)
#
#  - Adds a header label in column A, row 1 (matching the style of the
#    existing B1:E1 year/period headers) on each data sheet so the chart
#    series built off these tables gets a proper category/series header.
#  - Removes the (now redundant) bold/bordered header style from the
#    category cells in column A (rows 2..N), keeping their text, while
#    fixing missing accents on several labels.
#  - Sheet "Emissoes Totais (MtCO2eq)" also drops its trailing "Teto" row.
#  - Sheet "Custo Total (bilhões de R$)" gets a proper "2015" column
#    header (instead of "Custo") plus a new "Tipo Expansão" row header,
#    and updated cost figures.

$wb = $excel.ActiveWorkbook

function Set-HeaderCell {
    # Writes $text into $cellAddr and copies the cell formatting (only)
    # from $formatSourceAddr, so the new header cell matches the style
    # already used by the existing row/column headers.
    param($ws, [string]$cellAddr, [string]$text, [string]$formatSourceAddr)

    $ws.Range($cellAddr).Value = $text
    $ws.Range($formatSourceAddr).Copy() | Out-Null
    $ws.Range($cellAddr).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

function Clear-LabelStyle {
    # Optionally updates the cell text, then strips the bold/border
    # "header" formatting it used to carry.
    param($ws, [string]$cellAddr, [string]$text)

    if (-not [string]::IsNullOrEmpty($text)) {
        $ws.Range($cellAddr).Value = $text
    }
    $ws.Range($cellAddr).ClearFormats() | Out-Null
}

# --- Sheets 1-4: same A-column layout (Hidro, Gas Natural, Carvao, ...) ---
for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    Set-HeaderCell $ws "A1" "Fonte/Tecnologia" "B1"

    Clear-LabelStyle $ws "A2" ""                 # Hidro
    Clear-LabelStyle $ws "A3" "Gás Natural"       # Gas Natural -> Gás Natural
    Clear-LabelStyle $ws "A4" "Carvão"            # Carvao -> Carvão
    Clear-LabelStyle $ws "A5" ""                  # Nuclear
    Clear-LabelStyle $ws "A6" "Óleos Comb"        # Oleos Comb -> Óleos Comb
    Clear-LabelStyle $ws "A7" ""                  # Biomassa
    Clear-LabelStyle $ws "A8" "Eólica"            # Eolica -> Eólica
    Clear-LabelStyle $ws "A9" ""                  # Solar
    Clear-LabelStyle $ws "A10" ""                 # Outros
    Clear-LabelStyle $ws "A11" "Pot. Compl."      # Pot Compl -> Pot. Compl.
    Clear-LabelStyle $ws "A12" ""                 # GD
}

# --- Sheet 5: Emissoes Totais (MtCO2eq) -> P Medio / P Critico (+ drop Teto row) ---
$ws5 = $wb.Worksheets.Item(5)

Set-HeaderCell $ws5 "A1" "Período" "B1"

Clear-LabelStyle $ws5 "A2" "P.Médio"
Clear-LabelStyle $ws5 "A3" "P.Crítico"

$ws5.Rows(4).Delete() | Out-Null

# --- Sheet 6: Custo Total (bilhões de R$) ---
$ws6 = $wb.Worksheets.Item(6)

# B1 must stay a *text* "2015" header (not get auto-converted to a
# number), while keeping the existing bold/bordered header style - so
# force a text number format before assigning, then restore the header
# style (which the format change otherwise perturbs) from A2, which
# still carries the original, untouched header style at this point.
$ws6.Range("B1").NumberFormat = "@"
$ws6.Range("B1").Value = "2015"
$ws6.Range("A2").Copy() | Out-Null
$ws6.Range("B1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

Set-HeaderCell $ws6 "A1" "Tipo Expansão" "B1"

Clear-LabelStyle $ws6 "A2" "Expansão Centralizada"
$ws6.Range("B2").Value = 591

Clear-LabelStyle $ws6 "A3" "Expansão por GD"
$ws6.Range("B3").Value = 99
